$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct two existing data points (rows 69-70, Setting column E) ---
$ws.Range("E69").Value = 5230
$ws.Range("E70").Value = 6000

# --- Append new rows 623-638 with additional study data ---
$ws.Range("A623").Value = "Male"
$ws.Range("B623").Value = 20
$ws.Range("C623").Value = 60
$ws.Range("D623").Value = "maxFormantHz"
$ws.Range("E623").Value = 5230
$ws.Range("F623").Value = 45
$ws.Range("G623").Value = "10.1515/lingvan-2020-0051"
$ws.Range("A624").Value = "Female"
$ws.Range("B624").Value = 20
$ws.Range("C624").Value = 60
$ws.Range("D624").Value = "maxFormantHz"
$ws.Range("E624").Value = 6000
$ws.Range("F624").Value = 48
$ws.Range("G624").Value = "10.1515/lingvan-2020-0051"
$ws.Range("A625").Value = "Male"
$ws.Range("B625").Value = 10
$ws.Range("C625").Value = 12
$ws.Range("D625").Value = "maxFormantHz"
$ws.Range("E625").Value = 6300
$ws.Range("F625").Value = 27
$ws.Range("G625").Value = "10.1515/lingvan-2020-0051"
$ws.Range("A626").Value = "Female"
$ws.Range("B626").Value = 10
$ws.Range("C626").Value = 12
$ws.Range("D626").Value = "maxFormantHz"
$ws.Range("E626").Value = 6500
$ws.Range("F626").Value = 19
$ws.Range("G626").Value = "10.1515/lingvan-2020-0051"
$ws.Range("A627").Value = "Male"
$ws.Range("B627").Value = 20
$ws.Range("C627").Value = 60
$ws.Range("D627").Value = "nominalF1"
$ws.Range("E627").Formula = "=AVERAGE(342,768)"
$ws.Range("F627").Value = 45
$ws.Range("G627").Value = "10.1121/1.411872"
$ws.Range("A628").Value = "Male"
$ws.Range("B628").Value = 20
$ws.Range("C628").Value = 60
$ws.Range("D628").Value = "nominalF2"
$ws.Range("E628").Formula = "=AVERAGE(2322,910)"
$ws.Range("F628").Value = 45
$ws.Range("G628").Value = "10.1121/1.411872"
$ws.Range("A629").Value = "Male"
$ws.Range("B629").Value = 20
$ws.Range("C629").Value = 60
$ws.Range("D629").Value = "nominalF3"
$ws.Range("E629").Formula = "=AVERAGE(3000,1710)"
$ws.Range("F629").Value = 45
$ws.Range("G629").Value = "10.1121/1.411872"
$ws.Range("A630").Value = "Female"
$ws.Range("B630").Value = 20
$ws.Range("C630").Value = 60
$ws.Range("D630").Value = "nominalF1"
$ws.Range("E630").Formula = "=AVERAGE(936,437)"
$ws.Range("F630").Value = 48
$ws.Range("G630").Value = "10.1121/1.411872"
$ws.Range("A631").Value = "Female"
$ws.Range("B631").Value = 20
$ws.Range("C631").Value = 60
$ws.Range("D631").Value = "nominalF2"
$ws.Range("E631").Formula = "=AVERAGE(2761,1035)"
$ws.Range("F631").Value = 48
$ws.Range("G631").Value = "10.1121/1.411872"
$ws.Range("A632").Value = "Female"
$ws.Range("B632").Value = 20
$ws.Range("C632").Value = 60
$ws.Range("D632").Value = "nominalF3"
$ws.Range("E632").Formula = "=AVERAGE(3372,1929)"
$ws.Range("F632").Value = 48
$ws.Range("G632").Value = "10.1121/1.411872"
$ws.Range("A633").Value = "Male"
$ws.Range("B633").Value = 10
$ws.Range("C633").Value = 12
$ws.Range("D633").Value = "nominalF1"
$ws.Range("E633").Formula = "=AVERAGE(452,1002)"
$ws.Range("F633").Value = 27
$ws.Range("G633").Value = "10.1121/1.411872"
$ws.Range("A634").Value = "Male"
$ws.Range("B634").Value = 10
$ws.Range("C634").Value = 12
$ws.Range("D634").Value = "nominalF2"
$ws.Range("E634").Formula = "=AVERAGE(3081,1137)"
$ws.Range("F634").Value = 27
$ws.Range("G634").Value = "10.1121/1.411872"
$ws.Range("A635").Value = "Male"
$ws.Range("B635").Value = 10
$ws.Range("C635").Value = 12
$ws.Range("D635").Value = "nominalF3"
$ws.Range("E635").Formula = "=AVERAGE(3702,2950)"
$ws.Range("F635").Value = 27
$ws.Range("G635").Value = "10.1121/1.411872"
$ws.Range("A636").Value = "Female"
$ws.Range("B636").Value = 10
$ws.Range("C636").Value = 12
$ws.Range("D636").Value = "nominalF1"
$ws.Range("E636").Formula = "=AVERAGE(452,1002)"
$ws.Range("F636").Value = 19
$ws.Range("G636").Value = "10.1121/1.411872"
$ws.Range("A637").Value = "Female"
$ws.Range("B637").Value = 10
$ws.Range("C637").Value = 12
$ws.Range("D637").Value = "nominalF2"
$ws.Range("E637").Formula = "=AVERAGE(3081,1137)"
$ws.Range("F637").Value = 19
$ws.Range("G637").Value = "10.1121/1.411872"
$ws.Range("A638").Value = "Female"
$ws.Range("B638").Value = 10
$ws.Range("C638").Value = 12
$ws.Range("D638").Value = "nominalF3"
$ws.Range("E638").Formula = "=AVERAGE(3702,2950)"
$ws.Range("F638").Value = 19
$ws.Range("G638").Value = "10.1121/1.411872"

# --- Update selection / scroll position to match the final view state ---
$ws.Range("G630").Select()
